# Disponibilidad.xlsx update - automated availability check run
# Appends one more 14-row block (services Odoo..EZ Exporter) to the log,
# and corrects the timestamp on the previous block (D282:D295) to the
# precise value captured at save time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Fix up the previous block's timestamp (D282:D295): the value was
#    re-serialised with slightly higher precision on save.
# ---------------------------------------------------------------------
$ws.Range("D282:D295").Value2 = 44232.09276609954

# ---------------------------------------------------------------------
# 2) Append the new block: rows 296-309, same 14-service cycle used
#    throughout the sheet, all stamped with the same run timestamp.
# ---------------------------------------------------------------------
$names = @("Odoo", "Blackbox", "PowerBI", "Dropbox", "Odoo", "GEE", "UtilidadesOdoo", "Filtros Dashboard", "MapStore", "GeoServer", "Tomcat", "Shiny", "Github", "EZ Exporter")
$urls  = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$displayText = @(
    "https://www.dataintelligence-group.com/",
    "https://serviciodashboard.azurewebsites.net/",
    "https://powerbi.microsoft.com/es-es/",
    "https://www.dropbox.com/",
    "https://dataintelligence.store/",
    "https://app-data-i.users.earthengine.app/",
    "https://odooutil.azurewebsites.net/",
    "https://filtradordashboard.azurewebsites.net/",
    "https://ide.dataintelligence-group.com/mapstore/#/",
    "https://ide.dataintelligence-group.com/geoserver/web/?0",
    "https://ide.dataintelligence-group.com/",
    "https://rpubs.com/dataintelligence/",
    "https://github.com/Sud-Austral/",
    "https://ezexporter.highviewapps.com/exports/export-profile/"
)
$locations = @("", "", "", "", "", "", "", "", "/", "", "", "", "", "")
$status = "Disponible"
$stamp = 44232.113817718

$startRow = 296
for ($i = 0; $i -lt 14; $i++) {
    $r = $startRow + $i

    $ws.Cells.Item($r, 1).Value2 = $names[$i]
    $ws.Cells.Item($r, 3).Value2 = $status

    $d = $ws.Cells.Item($r, 4)
    $d.Value2 = $stamp
    $d.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $b = $ws.Cells.Item($r, 2)
    $b.Value2 = $displayText[$i]

    if ($locations[$i] -ne "") {
        $ws.Hyperlinks.Add($b, $urls[$i], $locations[$i])
    } else {
        $ws.Hyperlinks.Add($b, $urls[$i])
    }
    $b.Style = "Hyperlink"
}
